$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: rename "dimension"/refArea annotations to "measure" equivalents
$ws.Range("E2").Value = "iaest-measure:residencia-comarca-nombre"
$ws.Range("F2").Value = "iaest-measure:nacimiento-ccaa-nombre"
$ws.Range("G2").Value = "iaest-measure:sexo"
$ws.Range("I2").Value = "iaest-measure:residencia-provincia-nombre"
$ws.Range("J2").Value = "iaest-measure:residencia-ccaa-nombre"
$ws.Range("K2").Value = "iaest-measure:relacion-lugar-de-residencia-y-nacimiento"
$ws.Range("L2").Value = "iaest-measure:edad-grandes-grupos"

# Row 3: "dim" -> "medida" for the same columns
$ws.Range("E3").Value = "medida"
$ws.Range("F3").Value = "medida"
$ws.Range("G3").Value = "medida"
$ws.Range("I3").Value = "medida"
$ws.Range("J3").Value = "medida"
$ws.Range("K3").Value = "medida"
$ws.Range("L3").Value = "medida"

# Row 4: concept/URI type annotations replaced with "xsd:int"
$ws.Range("E4").Value = "xsd:int"
$ws.Range("F4").Value = "xsd:int"
$ws.Range("G4").Value = "xsd:int"
$ws.Range("I4").Value = "xsd:int"
$ws.Range("J4").Value = "xsd:int"
$ws.Range("K4").Value = "xsd:int"
$ws.Range("L4").Value = "xsd:int"

# Row 5: remove the mapping-file references that are no longer curated dimensions
$ws.Range("F5").Clear()
$ws.Range("G5").Clear()
$ws.Range("J5").Clear()
$ws.Range("K5").Clear()
$ws.Range("L5").Clear()
